$d = $word.ActiveDocument

# --- Helper: split a run by toggling font size away and back so the
# --- engine keeps the new piece as its own <w:r> instead of re-merging
# --- it into the neighbouring run that has identical formatting (the
# --- final size is restored, only the run boundary is what we want).
function Split-Run($rng) {
    $rng.Font.Size = 14
    $rng.Font.Size = 10.5
}

function Find-ParagraphByText($doc, $needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.TrimEnd() -eq $needle) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Change 1: paragraph "Products" gets a trailing space appended as its
# own run (same rPr: sz 21 / szCs 21 / lang en-US).
# ---------------------------------------------------------------------
$pProducts = Find-ParagraphByText $d 'Products'
if ($pProducts -ne $null) {
    $fullProducts = $pProducts.Range
    $oldEnd = $fullProducts.End
    $fullProducts.InsertAfter(' ')
    $rNew = $d.Range($oldEnd - 1, $oldEnd)
    Split-Run $rNew
}

# ---------------------------------------------------------------------
# Change 2: paragraph "Filter By Categories" -> "Filter By " + "Brand"
# (two runs).
# ---------------------------------------------------------------------
$pCategories = Find-ParagraphByText $d 'Filter By Categories'
if ($pCategories -ne $null) {
    $fullCategories = $pCategories.Range
    $start = $fullCategories.Start
    $rTail = $d.Range($start + 10, $start + 21)
    $rTail.Text = 'Brand'
    $rNewTail = $d.Range($start + 10, $start + 15)
    Split-Run $rNewTail
}

# ---------------------------------------------------------------------
# Change 3: paragraph "Filter By Rating" -> "Sort" + " By Rating"
# (two runs).
# ---------------------------------------------------------------------
$pRating = Find-ParagraphByText $d 'Filter By Rating'
if ($pRating -ne $null) {
    $fullRating = $pRating.Range
    $start = $fullRating.Start
    $rWord = $d.Range($start, $start + 6)
    $rWord.Text = 'Sort'
    $afterFull = $pRating.Range
    $rRest = $d.Range($start + 4, $afterFull.End)
    Split-Run $rRest
}
